$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Immunizations")

# Remove the "Status" column content from the Immunization Recommendations table
$ws.Range("C8:C10").ClearContents()

# A2: was bold+border -> now bold, no border
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Font.Bold = $true

# A3:A4: was gray bold wrap+border -> now gray bold wrap, no border
$r34 = $ws.Range("A3:A4")
$r34.Style = "Normal"
$r34.Font.Bold = $true
$r34.Font.Color = 6316128
$r34.WrapText = $true
$r34.VerticalAlignment = -4108

# A1 & A7: unchanged font (big bold 16pt), but xf index needs to shift down by one
$r17 = $ws.Range("A1")
$r17.Style = "Normal"
$r17.Font.Bold = $true
$r17.Font.Size = 16

$r7 = $ws.Range("A7")
$r7.Style = "Normal"
$r7.Font.Bold = $true
$r7.Font.Size = 16

# Update selection to reflect new active cell
$ws.Range("A11").Select()
